$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A9 value from -2 to -4
$ws.Range("A9").Value = -4

# Add a new row 12 with data (mirrors rows 8/10/11 pattern: date-label in column A)
$ws.Range("A12").Value = "23-02-2019"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 317
$ws.Range("J12").Value = 1
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 1
$ws.Range("P12").Value = 1
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = 2
$ws.Range("S12").Value = 2
$ws.Range("T12").Value = 0
